$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '92.328.29'
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("D3").Value = '3.090.98'
$ws.Range("E3").Value = '  -2.35%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '''238.45'
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").Value = '''608.64'
$ws.Range("E6").Value = '  -1.96%  '
$ws.Range("D7").Value = '''1.09'
$ws.Range("E7").Value = '  -2.38%  '
$ws.Range("E8").Value = '  +2.84%  '
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("D10").Value = '3.088.97'
$ws.Range("E10").Value = '  -2.30%  '
$ws.Range("D11").Value = '''0.746'
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("E12").Value = '  -1.78%  '
$ws.Range("E13").Value = '  -0.43%  '
$ws.Range("D14").Value = '92.559.11'
$ws.Range("E14").Value = '  +0.88%  '
$ws.Range("D15").Value = '''33.91'
$ws.Range("E15").Value = '  -4.68%  '
$ws.Range("E16").Value = '  -2.82%  '
$ws.Range("D17").Value = '3.685.47'
$ws.Range("E17").Value = '  -1.61%  '
$ws.Range("D18").Value = '3.100.95'
$ws.Range("E18").Value = '  -1.43%  '
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("E20").Value = '  -3.60%  '
$ws.Range("D21").Value = '''5.72'
$ws.Range("E21").Value = '  -3.79%  '
$ws.Range("D22").Value = '''443.30'
$ws.Range("E22").Value = '  -3.39%  '
$ws.Range("D23").Value = '''9.25'
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("D24").Value = '''0.0000195'
$ws.Range("E24").Value = '  -4.17%  '
$ws.Range("E25").Value = '  -5.89%  '
$ws.Range("D26").Value = '''86.10'
$ws.Range("E26").Value = '  -2.66%  '
$ws.Range("D27").Value = '''11.60'
$ws.Range("E27").Value = '  -3.95%  '
$ws.Range("D28").Value = '3.269.08'
$ws.Range("E28").Value = '  -1.52%  '
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("D30").Value = '''0.134'
$ws.Range("E30").Value = '  +5.48%  '
$ws.Range("D31").Value = '''0.227'
$ws.Range("E31").Value = '  -2.37%  '
$ws.Range("E32").Value = '  -2.38%  '
$ws.Range("D33").Value = '''9.05'
$ws.Range("E33").Value = '  -3.65%  '
$ws.Range("D34").Value = '''0.993'
$ws.Range("E34").Value = '  -0.71%  '
$ws.Range("D35").Value = '''7.92'
$ws.Range("E35").Value = '  +3.18%  '
$ws.Range("D36").Value = '''0.159'
$ws.Range("E36").Value = '  -7.14%  '
$ws.Range("D37").Value = '''25.84'
$ws.Range("E37").Value = '  -2.78%  '
$ws.Range("B38").Value = 'MantraDAO'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D38").Value = '''3.88'
$ws.Range("E38").Value = '  +1.81%  '
$ws.Range("B39").Value = 'PancakeSwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D39").Value = '''1.88'
$ws.Range("E39").Value = '  -2.83%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").Value = '''485.38'
$ws.Range("E40").Value = '  -5.39%  '
$ws.Range("B41").Value = 'WhiteBITCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D41").Value = '''23.85'
$ws.Range("E41").Value = '  +7.36%  '
$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").Value = '''1.28'
$ws.Range("E42").Value = '  -6.72%  '
$ws.Range("D43").Value = '''0.427'
$ws.Range("E43").Value = '  -5.40%  '
$ws.Range("D44").Value = '''3.28'
$ws.Range("E44").Value = '  -5.32%  '
$ws.Range("D46").Value = '''158.19'
$ws.Range("E46").Value = '  -0.99%  '
$ws.Range("E47").Value = '  -4.31%  '
$ws.Range("D48").Value = '''0.681'
$ws.Range("E48").Value = '  -4.99%  '
$ws.Range("E49").Value = '  -0.91%  '
$ws.Range("D50").Value = '''0.0328'
$ws.Range("E50").Value = '  +1.84%  '
$ws.Range("D51").Value = '''43.90'
$ws.Range("E51").Value = '  -0.57%  '
